$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. "April" -> "Linda" (both occurrences, same paragraph)
# ------------------------------------------------------------------
$d.Content.Find.Execute("April", $true, $false, $false, $false, $false, `
    $true, 1, $false, "Linda", 2) | Out-Null

# ------------------------------------------------------------------
# 2. Expand the "He must've dialed..." paragraph: insert a middle
#    clause and append the new "dull insectile ringing" material.
# ------------------------------------------------------------------
$d.Content.Find.Execute( `
    "redialed and set the receiver again to his ear. ", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "redialed checking the number magnet-stuck to the fridge even though he could easily recite it from memory and set the receiver again to his ear. The dull insectile ringing went on. Shelton imagined Linda shoving Jeff into the car, making plans to leave the baby with Jeff’s mother, already on her way up to the cabin. Pick up dammit. Pick up. I’m fine.  ", `
    2) | Out-Null

# ------------------------------------------------------------------
# 3. Insert the six brand-new paragraphs right after that paragraph.
# ------------------------------------------------------------------
$anchor = $d.Paragraphs(36)
$anchor.Range.InsertParagraphAfter()

$d.Paragraphs(37).Range.Text = "Should we go check?"
$d.Paragraphs(37).Range.InsertParagraphAfter()

$d.Paragraphs(38).Range.Text = "Yes. I’m feeling better. Why don’t we go check."
$d.Paragraphs(38).Range.InsertParagraphAfter()

$d.Paragraphs(39).Range.Text = "Shelton went to the bedroom to fortify his t-shirt and longjohns with jeans and a flannel shirt. Reaching into the closet for a shirt, Shelton saw the butt of the shotgun leaned into the far corner and concealed behind rows of his father’s clothes he thought he might someday wear. Shelton shut the closet doors."
$d.Paragraphs(39).Range.InsertParagraphAfter()

$d.Paragraphs(40).Range.Text = "The pilot waited beside the front door and examined Shelton as he reemerged into the living room. “It’s cold out there,” the pilot said reminding Shelton of his mother. As if he didn’t see the snow outside. As if he was going to head out in socks and no coat. Shelton snatched his Patagonia from a hook inside the pantry door and stepped into his Timberlands before wrapping the laces in bands around his ankles and tying them tightly."
$d.Paragraphs(40).Range.InsertParagraphAfter()

$d.Paragraphs(41).Range.Text = "“All set. Let’s go.”  The pilot opened the door inviting in a shrieking wind that pelted Shelton’s face with crystalline snow. "
$d.Paragraphs(41).Range.InsertParagraphAfter()

$d.Paragraphs(42).Range.Text = "They set off following the tracks the pilot had left. His footprints lay shallow and faint in a carpet of white having already been nearly erased by falling, drifting snows. The pilot never faltered though so Shelton stayed right behind him.   "

# ------------------------------------------------------------------
# 4. Tail whitespace paragraphs: the former lone-space paragraph
#    becomes a two-space paragraph, followed by a single-space
#    paragraph, then the pre-existing blank paragraph, then one more
#    new blank paragraph.
# ------------------------------------------------------------------
$d.Paragraphs(43).Range.Text = "  "
$d.Paragraphs(43).Range.InsertParagraphAfter()

$d.Paragraphs(44).Range.Text = " "

$d.Paragraphs(46).Range.InsertParagraphAfter()
